$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") holds text-formatted values (e.g. "72.372.15",
# "7.90") in the source data; force text format so Excel does not
# auto-convert numeric-looking strings to numbers (which would, e.g.,
# drop the trailing zero in "7.90").
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "72.372.15"
$ws.Range("E2").Value = "  -0.21%  "
$ws.Range("D3").Value = "2.639.03"
$ws.Range("E3").Value = "  -1.78%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "584.36"
$ws.Range("E5").Value = "  -2.78%  "
$ws.Range("D6").Value = "175.36"
$ws.Range("E6").Value = "  -1.41%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  -0.94%  "
$ws.Range("E9").Value = "  +1.18%  "
$ws.Range("D10").Value = "2.637.59"
$ws.Range("E10").Value = "  -1.80%  "
$ws.Range("E11").Value = "  +1.06%  "
$ws.Range("D12").Value = "0.358"
$ws.Range("E12").Value = "  +0.53%  "
$ws.Range("D13").Value = "4.93"
$ws.Range("E13").Value = "  -2.34%  "
$ws.Range("D14").Value = "3.121.75"
$ws.Range("E14").Value = "  -1.41%  "
$ws.Range("D15").Value = "0.0000186"
$ws.Range("E15").Value = "  +0.00%  "
$ws.Range("D16").Value = "72.197.85"
$ws.Range("E16").Value = "  -0.33%  "
$ws.Range("E17").Value = "  -2.06%  "
$ws.Range("D18").Value = "2.667.72"
$ws.Range("E18").Value = "  -0.43%  "
$ws.Range("E19").Value = "  +0.93%  "
$ws.Range("D20").Value = "7.90"
$ws.Range("E20").Value = "  -1.57%  "
$ws.Range("D21").Value = "376.07"
$ws.Range("E21").Value = "  +0.87%  "
$ws.Range("E22").Value = "  -1.46%  "
$ws.Range("D23").Value = "2.06"
$ws.Range("E23").Value = "  +0.95%  "
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("D25").Value = "71.06"
$ws.Range("E25").Value = "  -1.92%  "
$ws.Range("E26").Value = "  -2.33%  "
$ws.Range("D27").Value = "9.51"
$ws.Range("E27").Value = "  -3.71%  "
$ws.Range("D28").Value = "2.773.42"
$ws.Range("E28").Value = "  -1.23%  "
$ws.Range("E29").Value = "  +0.24%  "
$ws.Range("D30").Value = "0.0₃0951"
$ws.Range("E30").Value = "  +0.98%  "
$ws.Range("D31").Value = "7.96"
$ws.Range("D32").Value = "496.71"
$ws.Range("E32").Value = "  -3.19%  "
$ws.Range("E33").Value = "  -2.50%  "
$ws.Range("E34").Value = "  -1.76%  "
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("D36").Value = "160.25"
$ws.Range("E36").Value = "  -1.58%  "
$ws.Range("D37").Value = "19.19"
$ws.Range("E37").Value = "  -2.23%  "
$ws.Range("E38").Value = "  +5.49%  "
$ws.Range("E39").Value = "  -1.44%  "
$ws.Range("E40").Value = "  -1.68%  "
$ws.Range("E41").Value = "  -0.14%  "
$ws.Range("E42").Value = "  -5.26%  "
$ws.Range("E43").Value = "  -0.50%  "
$ws.Range("E44").Value = "  -2.80%  "
$ws.Range("E45").Value = "  -2.40%  "
$ws.Range("D46").Value = "39.01"
$ws.Range("E46").Value = "  -0.70%  "
$ws.Range("D47").Value = "151.51"
$ws.Range("E47").Value = "  -1.62%  "
$ws.Range("D48").Value = "3.65"
$ws.Range("E48").Value = "  -2.48%  "
$ws.Range("D49").Value = "0.544"
$ws.Range("E49").Value = "  -1.16%  "
$ws.Range("E50").Value = "  -3.70%  "
$ws.Range("D51").Value = "0.606"
$ws.Range("E51").Value = "  -0.22%  "
